$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "RP all" header to reflect the new LLNL comparisons
$ws.Range("J1").Value = "RP all (with new LLNL comparisons)"

# Updated validation numbers (latest V&V results)
$ws.Range("J3").Value = 1.27
$ws.Range("K3").Value = 0.59

# Move the active selection to J2, matching the saved view state
$ws.Range("J2").Select()
